$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 66800
$ws.Range("I19").Value = 200
$ws.Range("J19").Value = 100100
$ws.Range("K19").Value = 200
$ws.Range("L19").Value = 100100
$ws.Range("M19").Value = -25
$ws.Range("N19").Value = -100450
$ws.Range("H28").Value = 34503.766
$ws.Range("I28").Value = 44462.824
$ws.Range("K28").Value = 44462.824
$ws.Range("M28").Value = -43977.824
$ws.Range("H96").Value = 1515.2727
$ws.Range("I96").Value = 1870.8572
$ws.Range("K96").Value = 5612.571599999999
$ws.Range("M96").Value = -4239.571599999999
$ws.Range("H127").Value = 1313.0769
$ws.Range("I127").Value = 489.66666
$ws.Range("J127").Value = 2018.8572
$ws.Range("K127").Value = 1468.99998
$ws.Range("L127").Value = 6056.571599999999
$ws.Range("M127").Value = 3491.00002
$ws.Range("N127").Value = -15976.5716
$ws.Range("H132").Value = 12050664
$ws.Range("I132").Value = 13335630
$ws.Range("K132").Value = 40006890
$ws.Range("M132").Value = -40004360
$ws.Range("H138").Value = 2896.1406
$ws.Range("J138").Value = 3286.7437
$ws.Range("L138").Value = 9860.231100000001
$ws.Range("N138").Value = -20140.2311
$ws.Range("H141").Value = 2876.4473
$ws.Range("I141").Value = 2742.8857
$ws.Range("J141").Value = 4434.6665
$ws.Range("K141").Value = 8228.6571
$ws.Range("L141").Value = 13303.9995
$ws.Range("M141").Value = -3048.6571
$ws.Range("N141").Value = -23663.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15866.59
$ws.Range("I32").Value = 8446.207
$ws.Range("J32").Value = 37385.7
$ws.Range("K32").Value = 8446.207
$ws.Range("L32").Value = 37385.7
$ws.Range("M32").Value = -8159.207
$ws.Range("N32").Value = -37959.7
$ws.Range("H61").Value = 1783.4615
$ws.Range("I61").Value = 1580.3948
$ws.Range("K61").Value = 1580.3948
$ws.Range("M61").Value = -1368.3948
$ws.Range("H74").Value = 100005900
$ws.Range("I74").Value = 166671630
$ws.Range("J74").Value = 7323.5
$ws.Range("K74").Value = 166671630
$ws.Range("L74").Value = 7323.5
$ws.Range("M74").Value = -166670756
$ws.Range("N74").Value = -9071.5
$ws.Range("H77").Value = 100005900
$ws.Range("I77").Value = 166671630
$ws.Range("J77").Value = 7323.5
$ws.Range("K77").Value = 833358150
$ws.Range("L77").Value = 36617.5
$ws.Range("M77").Value = -833353782
$ws.Range("N77").Value = -45353.5
$ws.Range("H122").Value = 3308.4285
$ws.Range("I122").Value = 1844.4615
$ws.Range("J122").Value = 5687.375
$ws.Range("K122").Value = 5533.3845
$ws.Range("L122").Value = 17062.125
$ws.Range("M122").Value = -3083.3845
$ws.Range("N122").Value = -21962.125
$ws.Range("H135").Value = 27714.5
$ws.Range("J135").Value = 27714.5
$ws.Range("L135").Value = 27714.5
$ws.Range("N135").Value = -37854.5
$ws.Range("H136").Value = 1783.4615
$ws.Range("I136").Value = 1580.3948
$ws.Range("K136").Value = 4741.1844
$ws.Range("M136").Value = -2191.1844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2196.6936
$ws.Range("I134").Value = 1730.7142
$ws.Range("J134").Value = 6545.8335
$ws.Range("K134").Value = 5192.142599999999
$ws.Range("L134").Value = 19637.5005
$ws.Range("M134").Value = -2657.142599999999
$ws.Range("N134").Value = -24707.5005
$ws.Range("H135").Value = 49046.125
$ws.Range("J135").Value = 49046.125
$ws.Range("L135").Value = 49046.125
$ws.Range("N135").Value = -59186.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 1742.5
$ws.Range("I39").Value = 993
$ws.Range("K39").Value = 993
$ws.Range("M39").Value = -602
$ws.Range("H49").Value = 1742.5
$ws.Range("I49").Value = 993
$ws.Range("K49").Value = 993
$ws.Range("M49").Value = -811
$ws.Range("H94").Value = 878.5714
$ws.Range("I94").Value = 1191.3334
$ws.Range("J94").Value = 793.2727
$ws.Range("K94").Value = 1191.3334
$ws.Range("L94").Value = 793.2727
$ws.Range("M94").Value = -740.3334
$ws.Range("N94").Value = -1695.2727
$ws.Range("H99").Value = 333936.28
$ws.Range("I99").Value = 511347.72
$ws.Range("J99").Value = 8682
$ws.Range("K99").Value = 511347.72
$ws.Range("L99").Value = 8682
$ws.Range("M99").Value = -509849.72
$ws.Range("N99").Value = -11678
$ws.Range("H126").Value = 333936.28
$ws.Range("I126").Value = 511347.72
$ws.Range("J126").Value = 8682
$ws.Range("K126").Value = 1534043.16
$ws.Range("L126").Value = 26046
$ws.Range("M126").Value = -1531573.16
$ws.Range("N126").Value = -30986
$ws.Range("H132").Value = 420133.03
$ws.Range("I132").Value = 4454.6
$ws.Range("K132").Value = 13363.8
$ws.Range("M132").Value = -10833.8
$ws.Range("H134").Value = 3615.8367
$ws.Range("I134").Value = 3335.4688
$ws.Range("K134").Value = 10006.4064
$ws.Range("M134").Value = -7471.4064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 247.5
$ws.Range("J16").Value = 300.83334
$ws.Range("L16").Value = 902.5000200000001
$ws.Range("N16").Value = -1248.50002
$ws.Range("H132").Value = 2430.459
$ws.Range("I132").Value = 1611.9445
$ws.Range("K132").Value = 14507.5005
$ws.Range("M132").Value = -11977.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 58000
$ws.Range("J51").Value = 58000
$ws.Range("L51").Value = 58000
$ws.Range("N51").Value = -59018
$ws.Range("H80").Value = 3148.2
$ws.Range("I80").Value = 2005
$ws.Range("J80").Value = 3434
$ws.Range("K80").Value = 2005
$ws.Range("L80").Value = 3434
$ws.Range("M80").Value = -1007
$ws.Range("N80").Value = -5430
$ws.Range("H83").Value = 3148.2
$ws.Range("I83").Value = 2005
$ws.Range("J83").Value = 3434
$ws.Range("K83").Value = 10025
$ws.Range("L83").Value = 17170
$ws.Range("M83").Value = -5033
$ws.Range("N83").Value = -27154
$ws.Range("H122").Value = 267152.06
$ws.Range("I122").Value = 372546.72
$ws.Range("J122").Value = 3665.3333
$ws.Range("K122").Value = 1117640.16
$ws.Range("L122").Value = 10995.9999
$ws.Range("M122").Value = -1115190.16
$ws.Range("N122").Value = -15895.9999
$ws.Range("H132").Value = 3472.1052
$ws.Range("I132").Value = 3436.9033
$ws.Range("J132").Value = 3628
$ws.Range("K132").Value = 10310.7099
$ws.Range("L132").Value = 10884
$ws.Range("M132").Value = -7780.7099
$ws.Range("N132").Value = -15944

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 19015
$ws.Range("J57").Value = 19015
$ws.Range("L57").Value = 19015
$ws.Range("N57").Value = -20147
$ws.Range("H122").Value = 7440.0605
$ws.Range("I122").Value = 7263.524
$ws.Range("K122").Value = 21790.572
$ws.Range("M122").Value = -19340.572
$ws.Range("H132").Value = 16398246
$ws.Range("I132").Value = 23812658
$ws.Range("K132").Value = 71437974
$ws.Range("M132").Value = -71435444
$ws.Range("H133").Value = 94749.5
$ws.Range("J133").Value = 94749.5
$ws.Range("L133").Value = 94749.5
$ws.Range("N133").Value = -99809.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 69925
$ws.Range("J16").Value = 69925
$ws.Range("L16").Value = 69925
$ws.Range("N16").Value = -70509
$ws.Range("H126").Value = 1737.3334
$ws.Range("I126").Value = 1574.8
$ws.Range("K126").Value = 4724.4
$ws.Range("M126").Value = -2254.4
$ws.Range("H132").Value = 731406.9399999999
$ws.Range("I132").Value = 1118533.1
$ws.Range("J132").Value = 34579.8
$ws.Range("K132").Value = 3355599.3
$ws.Range("L132").Value = 103739.4
$ws.Range("M132").Value = -3353069.3
$ws.Range("N132").Value = -108799.4
$ws.Range("H136").Value = 4579.4546
$ws.Range("I136").Value = 4819.8975
$ws.Range("J136").Value = 3993.375
$ws.Range("K136").Value = 14459.6925
$ws.Range("L136").Value = 11980.125
$ws.Range("M136").Value = -11909.6925
$ws.Range("N136").Value = -17080.125

Write-Output "Applied all profit sheet updates"
